$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above row 223 (existing rows 223:233 shift down to 226:236)
$ws.Rows.Item(223).Resize(3).Insert()

# New row 223: Sandia, Extra, Region del Maule, fecha 2023-01-05 (44931)
$ws.Range("A223").Value = 5
$ws.Range("B223").Value = "Macroferia Regional de Talca"
$ws.Range("C223").Value = "Maule"
$ws.Range("D223").Value = 44931
$ws.Range("E223").Value = 7
$ws.Range("F223").Value = 100112028
$ws.Range("G223").Value = "Sandia"
$ws.Range("H223").Value = "Sin especificar"
$ws.Range("I223").Value = "Extra"
$ws.Range("J223").Value = 3000
$ws.Range("K223").Value = 2500
$ws.Range("L223").Value = 2500
$ws.Range("M223").Value = 2500
$ws.Range("N223").Value = "`$/unidad"
$ws.Range("O223").Value = "Región del Maule"
$ws.Range("P223").Value = 2500
$ws.Range("Q223").Value = 1
$ws.Range("R223").Value = "Hortaliza"

# New row 224: Sandia, Primera, Region del Maule, fecha 2023-01-05 (44931)
$ws.Range("A224").Value = 5
$ws.Range("B224").Value = "Macroferia Regional de Talca"
$ws.Range("C224").Value = "Maule"
$ws.Range("D224").Value = 44931
$ws.Range("E224").Value = 7
$ws.Range("F224").Value = 100112028
$ws.Range("G224").Value = "Sandia"
$ws.Range("H224").Value = "Sin especificar"
$ws.Range("I224").Value = "Primera"
$ws.Range("J224").Value = 2000
$ws.Range("K224").Value = 2000
$ws.Range("L224").Value = 2000
$ws.Range("M224").Value = 2000
$ws.Range("N224").Value = "`$/unidad"
$ws.Range("O224").Value = "Región del Maule"
$ws.Range("P224").Value = 2000
$ws.Range("Q224").Value = 1
$ws.Range("R224").Value = "Hortaliza"

# New row 225: Sandia, Segunda, Region del Maule, fecha 2023-01-05 (44931)
$ws.Range("A225").Value = 5
$ws.Range("B225").Value = "Macroferia Regional de Talca"
$ws.Range("C225").Value = "Maule"
$ws.Range("D225").Value = 44931
$ws.Range("E225").Value = 7
$ws.Range("F225").Value = 100112028
$ws.Range("G225").Value = "Sandia"
$ws.Range("H225").Value = "Sin especificar"
$ws.Range("I225").Value = "Segunda"
$ws.Range("J225").Value = 2000
$ws.Range("K225").Value = 1500
$ws.Range("L225").Value = 1500
$ws.Range("M225").Value = 1500
$ws.Range("N225").Value = "`$/unidad"
$ws.Range("O225").Value = "Región del Maule"
$ws.Range("P225").Value = 1500
$ws.Range("Q225").Value = 1
$ws.Range("R225").Value = "Hortaliza"
